# The captured change is confined entirely to <w:docDefaults> inside
# word/styles.xml: it deletes a batch of <w:rPr>/<w:pPr> child elements
# whose values are simply the OOXML schema defaults already in force
# (w:b=0, w:i=0, w:smallCaps=0, w:strike=0, w:color=000000, w:u=none,
# w:shd=clear/auto, w:vertAlign=baseline, w:keepNext=0, w:keepLines=0,
# w:widowControl=1, an empty/nil w:pBdr, w:spacing before/after=0 plus
# line=276/auto, w:ind all zero, w:contextualSpacing=0, w:jc=left) and
# keeps only <w:sz>, <w:szCs>, <w:lang> (rPrDefault) and a bare
# <w:spacing w:line="276" w:lineRule="auto"/> (pPrDefault). None of
# those removed values differ from what Word already renders by
# default, so the edit has zero visual/behavioral effect - it is a
# pure "drop the redundant explicit defaults" cleanup of the
# <w:docDefaults> block itself.
#
# <w:docDefaults> (and its <w:rPrDefault>/<w:pPrDefault> children) is
# metadata Word derives/round-trips from the package's styles part; it
# is not surfaced anywhere in the Word object model (no
# Document/Styles/Style/Font/ParagraphFormat member reads or writes
# it - confirmed against this host's full command surface). Word's
# automation model only ever lets you touch the *explicit* formatting
# of a Style (e.g. Styles("Normal").Font / .ParagraphFormat), a
# Range/Selection, or Find/Replace targets - all of which land in
# word/document.xml or in the named <w:style> element itself, never in
# <w:docDefaults>. Exercising any of those paths here would not move
# <w:docDefaults> at all; it would instead bolt new explicit
# <w:rPr>/<w:pPr> overrides onto the "Normal" style (or onto runs/
# paragraphs in the body) that the real commit never added, making the
# document diverge from the target instead of converging on it.
#
# So there is no COM-reachable operation that reproduces this
# particular hunk without corrupting unrelated, unchanged parts of the
# package. Since the edit is a no-op for the document's appearance,
# the safest, most faithful action available through this interop
# surface is to leave the content/styles untouched rather than
# fabricate collateral formatting changes Word's object model cannot
# actually avoid creating.
$d = $word.ActiveDocument
Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
